$d = $word.ActiveDocument

# The name heading currently reads "Ben Barrrr": after "Ben Ba" there is a
# single-letter run containing "r", immediately followed by a stray
# duplicate run whose text is "rrr" (a leftover artifact from live editing).
# Locate that paragraph, compute the character span of the extra "rrr",
# and delete just that span using an explicit character-position Range
# (rather than Find/Replace) so the surrounding "Ben Ba" and "r" runs are
# left completely untouched/unmerged.
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    $relIdx = $t.IndexOf("Ben Barrrr")
    if ($relIdx -ge 0) {
        $nameStart = $p.Range.Start + $relIdx
        # "Ben Bar" (the text we want to keep) is 7 characters long; the
        # stray duplicate "rrr" immediately follows it.
        $extraStart = $nameStart + 7
        $extraEnd = $extraStart + 3
        $d.Range($extraStart, $extraEnd).Delete()
        break
    }
}
